$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.374.54"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "'2.064.94"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'233.92"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'56.74"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "'2.368.13"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "'14.58"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "'20.65"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "'0.778"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("D17").Value = "'2.064.19"
$ws.Range("E17").Value = "  -2.34%  "
$ws.Range("D18").Value = "'37.306.25"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("D20").Value = "'69.45"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'226.31"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("D26").Value = "'166.62"
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("D27").Value = "'8.78"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("E28").Value = "  +4.69%  "
$ws.Range("D29").Value = "'19.05"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").Value = "'0.127"
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").Value = "'4.56"
$ws.Range("E34").Value = "  +3.82%  "
$ws.Range("D35").Value = "'2.48"
$ws.Range("E35").Value = "  -2.43%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "'3.24"
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("D40").Value = "'2.95"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("D42").Value = "'1.463.91"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "'4.33"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.0934"
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'1.17"
$ws.Range("E45").Value = "  +3.71%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0212"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'15.04"
$ws.Range("E48").Value = "  -7.17%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'7.13"
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.95"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").Value = "'2.254.26"
$ws.Range("E51").Value = "  -0.69%  "
